$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.733.77'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('D3').Value = '2.553.90'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '528.62'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('D6').Value = '134.68'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('D8').Value = '0.567'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('D9').Value = '2.546.46'
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('D10').Value = '0.0992'
$ws.Range('E10').Value = '  +0.33%  '
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('D12').Value = '5.21'
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('D13').Value = '0.335'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').Value = '2.998.63'
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('D15').Value = '58.720.53'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').Value = '22.45'
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '0.0000137'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = '2.545.15'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').Value = '10.78'
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Value = '324.79'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('D21').Value = '4.21'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('D22').Value = '6.18'
$ws.Range('E22').Value = '  +7.16%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = '65.34'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('D25').Value = '0.413'
$ws.Range('E25').Value = '  +0.07%  '
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = '7.47'
$ws.Range('E28').Value = '  +1.02%  '
$ws.Range('D29').Value = '0.0₃0765'
$ws.Range('E29').Value = '  +1.73%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '1.22'
$ws.Range('E30').Value = '  +2.41%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.74'
$ws.Range('E31').Value = '  +2.33%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').Value = '6.41'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '167.97'
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = '0.995'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('D36').Value = '18.36'
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('D37').Value = '1.28'
$ws.Range('E37').Value = '  -2.95%  '
$ws.Range('D38').Value = '4.00'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('D39').Value = '1.51'
$ws.Range('E39').Value = '  +2.16%  '
$ws.Range('D40').Value = '36.65'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').Value = '0.794'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D42').Value = '283.16'
$ws.Range('E42').Value = '  +2.38%  '
$ws.Range('D43').Value = '3.50'
$ws.Range('E43').Value = '  +1.26%  '
$ws.Range('D44').Value = '5.13'
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = '0.606'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '130.32'
$ws.Range('E46').Value = '  +4.69%  '
$ws.Range('D47').Value = '0.0920'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').Value = '0.0507'
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range('D49').Value = '17.99'
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('D50').Value = '0.0217'
$ws.Range('E50').Value = '  +1.18%  '
$ws.Range('D51').Value = '17.36'
$ws.Range('E51').Value = '  +1.06%  '
